$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes (shared strings) ---
# Row 4 (BR50043): Reason + Error Message now reference the shortened "both URLs failed" text
$ws.Range("E4").Value2 = "both URLs failed"
$ws.Range("F4").Value2 = "Failed to download PDF: both URLs failed"

# Row 5 (BR50044): Error Message reason changed to a different failure description
$ws.Range("F5").Value2 = "Failed to download PDF: first URL failed and no fallback"

# Row 11 (BR50050): Reason cleared (was re-using the old "download failed - both URLs tried" text)
$ws.Range("E11").Value2 = ""

# --- Column width changes ---
$ws.Columns.Item(5).ColumnWidth = 24.333333
$ws.Columns.Item(6).ColumnWidth = 49.166667
